# Apply the "Update countries & provincias Spain" data refresh to the Pais sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "last updated" timestamp shown in A1.
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 12:58"

# Refresh per-country COVID figures (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes). The data is kept sorted by
# total cases, so a couple of rows now hold a different country than before:
#   row 46/47  -> Rumania overtakes Portugal
#   row 79/80  -> Estado de Palestina overtakes Bosnia y Herzegovina

$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 1701532
$ws.Range("C6").Value = 4478
$ws.Range("D6").Value = 1096898
$ws.Range("E6").Value = 568047
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 36
$ws.Range("H6").Value = 36587

$ws.Range("A13").Value = "Iran"
$ws.Range("B13").Value = 306752
$ws.Range("C13").Value = 2548
$ws.Range("D13").Value = 265830
$ws.Range("E13").Value = 23940
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 216
$ws.Range("H13").Value = 16982

$ws.Range("A43").Value = "Emiratos Arabes Unidos"
$ws.Range("B43").Value = 60760
$ws.Range("C43").Value = 254
$ws.Range("D43").Value = 54255
$ws.Range("E43").Value = 6154
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 351

$ws.Range("A46").Value = "Rumania"
$ws.Range("B46").Value = 52111
$ws.Range("C46").Value = 1225
$ws.Range("D46").Value = 27346
$ws.Range("E46").Value = 22386
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 36
$ws.Range("H46").Value = 2379

$ws.Range("A47").Value = "Portugal"
$ws.Range("B47").Value = 51072
$ws.Range("C47").Value = 0
$ws.Range("D47").Value = 36483
$ws.Range("E47").Value = 12854
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 1735

$ws.Range("A57").Value = "Suiza"
$ws.Range("B57").Value = 35412
$ws.Range("C57").Value = 180
$ws.Range("D57").Value = 31100
$ws.Range("E57").Value = 2331
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 1981

$ws.Range("A74").Value = "El Salvador"
$ws.Range("B74").Value = 17050
$ws.Range("C74").Value = 418
$ws.Range("D74").Value = 8495
$ws.Range("E74").Value = 8096
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 11
$ws.Range("H74").Value = 459

$ws.Range("A79").Value = "Estado de Palestina"
$ws.Range("B79").Value = 12160
$ws.Range("C79").Value = 323
$ws.Range("D79").Value = 5324
$ws.Range("E79").Value = 6754
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 82

$ws.Range("A80").Value = "Bosnia y Herzegovina"
$ws.Range("B80").Value = 11876
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 5959
$ws.Range("E80").Value = 5578
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 339

$ws.Range("A85").Value = "Senegal"
$ws.Range("B85").Value = 10284
$ws.Range("C85").Value = 52
$ws.Range("D85").Value = 6822
$ws.Range("E85").Value = 3253
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 4
$ws.Range("H85").Value = 209

$ws.Range("A97").Value = "Zambia"
$ws.Range("B97").Value = 6228
$ws.Range("C97").Value = 265
$ws.Range("D97").Value = 4130
$ws.Range("E97").Value = 1933
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 14
$ws.Range("H97").Value = 165

$ws.Range("A154").Value = "Malta"
$ws.Range("B154").Value = 845
$ws.Range("C154").Value = 21
$ws.Range("D154").Value = 665
$ws.Range("E154").Value = 171
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 9
